$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5120
$ws.Range("I12").Value = 4816.6665
$ws.Range("J12").Value = 5666
$ws.Range("K12").Value = 4816.6665
$ws.Range("L12").Value = 5666
$ws.Range("M12").Value = -4646.6665
$ws.Range("N12").Value = -6006

$ws.Range("H70").Value = 251874.75
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 251874.75
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H76").Value = 111127740
$ws.Range("I76").Value = 10670.857
$ws.Range("J76").Value = 500037500
$ws.Range("K76").Value = 10670.857
$ws.Range("L76").Value = 500037500
$ws.Range("M76").Value = -10355.857
$ws.Range("N76").Value = -500038130

$ws.Range("H79").Value = 111127740
$ws.Range("I79").Value = 10670.857
$ws.Range("J79").Value = 500037500
$ws.Range("K79").Value = 10670.857
$ws.Range("L79").Value = 500037500
$ws.Range("M79").Value = -9578.857
$ws.Range("N79").Value = -500039684

$ws.Range("H100").Value = 5391.3335
$ws.Range("I100").Value = 1851.2
$ws.Range("K100").Value = 1851.2
$ws.Range("M100").Value = -1310.2

$ws.Range("H103").Value = 428.96667
$ws.Range("I103").Value = 367.3
$ws.Range("J103").Value = 552.3
$ws.Range("K103").Value = 1101.9
$ws.Range("L103").Value = 1656.9
$ws.Range("M103").Value = -515.9000000000001
$ws.Range("N103").Value = -2828.9

$ws.Range("H127").Value = 15993.5
$ws.Range("I127").Value = 17820.5
$ws.Range("K127").Value = 53461.5
$ws.Range("M127").Value = -48501.5

$ws.Range("H132").Value = 2858.4783
$ws.Range("I132").Value = 1574.7059
$ws.Range("K132").Value = 4724.1177
$ws.Range("M132").Value = -2194.1177

$ws.Range("H135").Value = 1416.75
$ws.Range("I135").Value = 695.92
$ws.Range("J135").Value = 3991.1428
$ws.Range("K135").Value = 6263.28
$ws.Range("L135").Value = 35920.2852
$ws.Range("M135").Value = -3728.28
$ws.Range("N135").Value = -40990.2852

$ws.Range("H137").Value = 4376.7754
$ws.Range("I137").Value = 2903.4583
$ws.Range("K137").Value = 8710.374899999999
$ws.Range("M137").Value = -6160.374899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 45524.5
$ws.Range("J44").Value = 45524.5
$ws.Range("L44").Value = 45524.5
$ws.Range("N44").Value = -46500.5

$ws.Range("H55").Value = 73276.5
$ws.Range("J55").Value = 73276.5
$ws.Range("L55").Value = 73276.5
$ws.Range("N55").Value = -73906.5

$ws.Range("H61").Value = 2315.9575
$ws.Range("I61").Value = 1581.8292
$ws.Range("K61").Value = 1581.8292
$ws.Range("M61").Value = -1369.8292

$ws.Range("H74").Value = 1386.683
$ws.Range("I74").Value = 1185.742
$ws.Range("J74").Value = 2009.6
$ws.Range("K74").Value = 1185.742
$ws.Range("L74").Value = 2009.6
$ws.Range("M74").Value = -311.742
$ws.Range("N74").Value = -3757.6

$ws.Range("H77").Value = 1386.683
$ws.Range("I77").Value = 1185.742
$ws.Range("J77").Value = 2009.6
$ws.Range("K77").Value = 5928.71
$ws.Range("L77").Value = 10048
$ws.Range("M77").Value = -1560.71
$ws.Range("N77").Value = -18784

$ws.Range("H102").Value = 1656.1333
$ws.Range("I102").Value = 1663.9656
$ws.Range("K102").Value = 1663.9656
$ws.Range("M102").Value = -41.96559999999999

$ws.Range("H122").Value = 3086.8298
$ws.Range("I122").Value = 1981
$ws.Range("K122").Value = 5943
$ws.Range("M122").Value = -3493

$ws.Range("H136").Value = 2315.9575
$ws.Range("I136").Value = 1581.8292
$ws.Range("K136").Value = 4745.487599999999
$ws.Range("M136").Value = -2195.487599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2852.8333
$ws.Range("I20").Value = 2310.2222
$ws.Range("J20").Value = 3395.4443
$ws.Range("K20").Value = 2310.2222
$ws.Range("L20").Value = 3395.4443
$ws.Range("M20").Value = -2063.2222
$ws.Range("N20").Value = -3889.4443

$ws.Range("H86").Value = 1312.6
$ws.Range("J86").Value = 415
$ws.Range("L86").Value = 415
$ws.Range("N86").Value = -2661

$ws.Range("H89").Value = 1312.6
$ws.Range("J89").Value = 415
$ws.Range("L89").Value = 2075
$ws.Range("N89").Value = -13307

$ws.Range("H94").Value = 962.8889
$ws.Range("I94").Value = 962.8889
$ws.Range("K94").Value = 962.8889
$ws.Range("M94").Value = -511.8889

$ws.Range("H105").Value = 75801.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 628
$ws.Range("I22").Value = 285
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 285
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 65
$ws.Range("N22").Value = -2700

$ws.Range("H59").Value = 87557.5
$ws.Range("J59").Value = 87557.5
$ws.Range("L59").Value = 87557.5
$ws.Range("N59").Value = -89847.5

$ws.Range("H122").Value = 2842.5334
$ws.Range("I122").Value = 1451.25
$ws.Range("K122").Value = 4353.75
$ws.Range("M122").Value = -1903.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 72500
$ws.Range("J42").Value = 72500
$ws.Range("L42").Value = 72500
$ws.Range("N42").Value = -73470

$ws.Range("H115").Value = 72500
$ws.Range("J115").Value = 72500
$ws.Range("L115").Value = 72500
$ws.Range("N115").Value = -74850

$ws.Range("H132").Value = 531274
$ws.Range("I132").Value = 670980.5600000001
$ws.Range("J132").Value = 7374.5
$ws.Range("K132").Value = 2012941.68
$ws.Range("L132").Value = 22123.5
$ws.Range("M132").Value = -2010411.68
$ws.Range("N132").Value = -27183.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3206.8823
$ws.Range("J46").Value = 3949.5
$ws.Range("L46").Value = 3949.5
$ws.Range("N46").Value = -4325.5

$ws.Range("H122").Value = 1169556.2
$ws.Range("I122").Value = 836444.75
$ws.Range("K122").Value = 2509334.25
$ws.Range("M122").Value = -2506884.25

$ws.Range("H132").Value = 4080.85
$ws.Range("I132").Value = 2494.1
$ws.Range("J132").Value = 5667.6
$ws.Range("K132").Value = 7482.299999999999
$ws.Range("L132").Value = 17002.8
$ws.Range("M132").Value = -4952.299999999999
$ws.Range("N132").Value = -22062.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 78709.5
$ws.Range("J124").Value = 78709.5
$ws.Range("L124").Value = 78709.5
$ws.Range("N124").Value = -88529.5

$ws.Range("H129").Value = 49916.668
$ws.Range("J129").Value = 49916.668
$ws.Range("L129").Value = 49916.668
$ws.Range("N129").Value = -59916.668
